# Add a new "PWS" (solar) capacity-factor block, mirroring the existing
# "Solar"/"Wind" CF columns (H:I), into columns K:M, plus a second
# Correct-CF/Eff/EFF_TIMES projection table at K22:M25 (mirrors K16:M20,
# but driven off the new PWS total in L12). Also clears the now-unused
# literal numbers that used to sit in B15:B22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row 3: new "PWS" / "CF" columns, bold like the existing headers.
# ---------------------------------------------------------------------
$ws.Range("K3").Value = "PWS"
$ws.Range("K3").Font.Bold = $true
$ws.Range("L3").Value = "CF"
$ws.Range("L3").Font.Bold = $true

# ---------------------------------------------------------------------
# K4:K11 - PWS on/off flags (copy of the H column pattern).
# L4:L10 - shared formula K*C; L11 stands out (K11*F11, matching source).
# M4:M11 - labels, same text (and same shared-string slot) as column A.
# ---------------------------------------------------------------------
$kFlags = @(1, 1, 1, 0, 1, 0, 1, 0)
for ($i = 0; $i -lt 8; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 11).Value = $kFlags[$i]
}

$ws.Range("L4").Formula = "=K4*C4"
$ws.Range("L5:L10").Formula = "=K5*C5"
$ws.Range("L11").Formula = "=K11*F11"

$labels = @("SDW", "SNW", "WDW", "WNW", "SDN", "SNN", "WDN", "WNN")
for ($i = 0; $i -lt 8; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 13).Value = $labels[$i]
    $ws.Cells.Item($row, 13).Font.Bold = $true
}

# Row 12 total for the new PWS CF column.
$ws.Range("L12").Formula = "=SUM(L4:L11)"

# ---------------------------------------------------------------------
# The old ad-hoc numbers in B15:B22 are no longer used - clear them but
# keep their existing number format.
# ---------------------------------------------------------------------
foreach ($r in 15..22) {
    $ws.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------
# Second correct-CF / efficiency / EFF_TIMES table at K22:M25, mirroring
# K16:M20 but referencing the new PWS total ($L$12) and a PWS row label.
# ---------------------------------------------------------------------
$ws.Range("K22").Value = "Correct CF"
$ws.Range("K22").Font.Bold = $true
$ws.Range("L22").Value = "Eff"
$ws.Range("L22").Font.Bold = $true
$ws.Range("M22").Value = "EFF_TIMES"
$ws.Range("M22").Font.Bold = $true

$ws.Range("K23").Value = 0.21
$ws.Range("L23").Value = 0.15
$ws.Range("M23").Formula = "=K23/`$L`$12"

$ws.Range("I24").Value = "PWS"
$ws.Range("I24").Font.Bold = $true
$ws.Range("K24").Formula = "=K23*(1+(L24-L23))"
$ws.Range("L24").Value = 0.25
$ws.Range("M24:M25").Formula = "=K24/`$L`$12"

$ws.Range("K25").Formula = "=K24*(1+(L25-L24))"
$ws.Range("L25").Value = 0.1

# Apply the 0.00 number format to L24 last - doing this earlier causes the
# engine to leak the style onto the next-written formula cell (K25).
$ws.Range("L24").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# Approximate the saved view state (best-effort; exact anchor cell
# within a multi-cell selection isn't independently scriptable).
# ---------------------------------------------------------------------
$ws.Range("M23:M25").Select()
